$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(7)

# --- Column width adjustments (col B / X / Y) ---
# Target widths (from source XML) are on a non-Excel quantization grid; these
# inputs are chosen so that, after COM's internal round(w*7)+5 pixel quantization,
# the stored width lands on the closest achievable 1/7-character-unit value.
$ws.Columns.Item(2).ColumnWidth = 7.857142857142857
$ws.Columns.Item(2).Hidden = $true
$ws.Columns.Item(24).ColumnWidth = 46.42857142857143
$ws.Columns.Item(25).ColumnWidth = 70.57142857142857

# --- Row rewrites: each row below stored its 6 meaningful fields packed into a
# single weird array-string in column B. Unpack them into the normal per-column
# layout (same columns used by every other data row: W/X/Y text + blank Z/AA/AB,
# blank C..V with M/U kept in the date style used throughout the sheet).

# Row 7
$ws.Cells.Item(7, 2).Value = ""
$ws.Cells.Item(7, 2).Style = "Normal"
$ws.Cells.Item(7, 3).Style = "Normal"
$ws.Cells.Item(7, 4).Style = "Normal"
$ws.Cells.Item(7, 5).Style = "Normal"
$ws.Cells.Item(7, 6).Style = "Normal"
$ws.Cells.Item(7, 7).Style = "Normal"
$ws.Cells.Item(7, 8).Style = "Normal"
$ws.Cells.Item(7, 9).Style = "Normal"
$ws.Cells.Item(7, 10).Style = "Normal"
$ws.Cells.Item(7, 11).Style = "Normal"
$ws.Cells.Item(7, 12).Style = "Normal"
$ws.Cells.Item(7, 14).Style = "Normal"
$ws.Cells.Item(7, 15).Style = "Normal"
$ws.Cells.Item(7, 16).Style = "Normal"
$ws.Cells.Item(7, 17).Style = "Normal"
$ws.Cells.Item(7, 18).Style = "Normal"
$ws.Cells.Item(7, 19).Style = "Normal"
$ws.Cells.Item(7, 20).Style = "Normal"
$ws.Cells.Item(7, 22).Style = "Normal"
$ws.Cells.Item(7, 13).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(7, 21).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(7, 23).Value = 'PANELS-(PANELS)'
$ws.Cells.Item(7, 24).Value = 'Left side panel-(Left side panel)'
$ws.Cells.Item(7, 25).Value = 'LEFT SIDE PANEL CUT 6" (GUM FITTINGS)'
$ws.Cells.Item(7, 26).Value = ""
$ws.Cells.Item(7, 26).Style = "Normal"
$ws.Cells.Item(7, 27).Value = ""
$ws.Cells.Item(7, 27).Style = "Normal"
$ws.Cells.Item(7, 28).Value = ""
$ws.Cells.Item(7, 28).Style = "Normal"

# Row 9
$ws.Cells.Item(9, 2).Value = ""
$ws.Cells.Item(9, 2).Style = "Normal"
$ws.Cells.Item(9, 3).Style = "Normal"
$ws.Cells.Item(9, 4).Style = "Normal"
$ws.Cells.Item(9, 5).Style = "Normal"
$ws.Cells.Item(9, 6).Style = "Normal"
$ws.Cells.Item(9, 7).Style = "Normal"
$ws.Cells.Item(9, 8).Style = "Normal"
$ws.Cells.Item(9, 9).Style = "Normal"
$ws.Cells.Item(9, 10).Style = "Normal"
$ws.Cells.Item(9, 11).Style = "Normal"
$ws.Cells.Item(9, 12).Style = "Normal"
$ws.Cells.Item(9, 14).Style = "Normal"
$ws.Cells.Item(9, 15).Style = "Normal"
$ws.Cells.Item(9, 16).Style = "Normal"
$ws.Cells.Item(9, 17).Style = "Normal"
$ws.Cells.Item(9, 18).Style = "Normal"
$ws.Cells.Item(9, 19).Style = "Normal"
$ws.Cells.Item(9, 20).Style = "Normal"
$ws.Cells.Item(9, 22).Style = "Normal"
$ws.Cells.Item(9, 13).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(9, 21).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(9, 23).Value = 'DOORS-(D)'
$ws.Cells.Item(9, 24).Value = 'Door stiffeners hinges side edge.-(Door stiffeners hinges side edge.)'
$ws.Cells.Item(9, 25).Value = 'LEFT DOOR LOCK BAR BENT 01 PCS'
$ws.Cells.Item(9, 26).Value = ""
$ws.Cells.Item(9, 26).Style = "Normal"
$ws.Cells.Item(9, 27).Value = ""
$ws.Cells.Item(9, 27).Style = "Normal"
$ws.Cells.Item(9, 28).Value = ""
$ws.Cells.Item(9, 28).Style = "Normal"

# Row 12
$ws.Cells.Item(12, 2).Value = ""
$ws.Cells.Item(12, 2).Style = "Normal"
$ws.Cells.Item(12, 3).Style = "Normal"
$ws.Cells.Item(12, 4).Style = "Normal"
$ws.Cells.Item(12, 5).Style = "Normal"
$ws.Cells.Item(12, 6).Style = "Normal"
$ws.Cells.Item(12, 7).Style = "Normal"
$ws.Cells.Item(12, 8).Style = "Normal"
$ws.Cells.Item(12, 9).Style = "Normal"
$ws.Cells.Item(12, 10).Style = "Normal"
$ws.Cells.Item(12, 11).Style = "Normal"
$ws.Cells.Item(12, 12).Style = "Normal"
$ws.Cells.Item(12, 14).Style = "Normal"
$ws.Cells.Item(12, 15).Style = "Normal"
$ws.Cells.Item(12, 16).Style = "Normal"
$ws.Cells.Item(12, 17).Style = "Normal"
$ws.Cells.Item(12, 18).Style = "Normal"
$ws.Cells.Item(12, 19).Style = "Normal"
$ws.Cells.Item(12, 20).Style = "Normal"
$ws.Cells.Item(12, 22).Style = "Normal"
$ws.Cells.Item(12, 13).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(12, 21).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(12, 23).Value = 'DOORS-(D)'
$ws.Cells.Item(12, 24).Value = 'Door stiffeners hinges side edge.-(Door stiffeners hinges side edge.)'
$ws.Cells.Item(12, 25).Value = 'RIGHT DOOR BOTTOM GASKET LOOSE FULLY'
$ws.Cells.Item(12, 26).Value = ""
$ws.Cells.Item(12, 26).Style = "Normal"
$ws.Cells.Item(12, 27).Value = ""
$ws.Cells.Item(12, 27).Style = "Normal"
$ws.Cells.Item(12, 28).Value = ""
$ws.Cells.Item(12, 28).Style = "Normal"

# Row 15
$ws.Cells.Item(15, 2).Value = ""
$ws.Cells.Item(15, 2).Style = "Normal"
$ws.Cells.Item(15, 3).Style = "Normal"
$ws.Cells.Item(15, 4).Style = "Normal"
$ws.Cells.Item(15, 5).Style = "Normal"
$ws.Cells.Item(15, 6).Style = "Normal"
$ws.Cells.Item(15, 7).Style = "Normal"
$ws.Cells.Item(15, 8).Style = "Normal"
$ws.Cells.Item(15, 9).Style = "Normal"
$ws.Cells.Item(15, 10).Style = "Normal"
$ws.Cells.Item(15, 11).Style = "Normal"
$ws.Cells.Item(15, 12).Style = "Normal"
$ws.Cells.Item(15, 14).Style = "Normal"
$ws.Cells.Item(15, 15).Style = "Normal"
$ws.Cells.Item(15, 16).Style = "Normal"
$ws.Cells.Item(15, 17).Style = "Normal"
$ws.Cells.Item(15, 18).Style = "Normal"
$ws.Cells.Item(15, 19).Style = "Normal"
$ws.Cells.Item(15, 20).Style = "Normal"
$ws.Cells.Item(15, 22).Style = "Normal"
$ws.Cells.Item(15, 13).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(15, 21).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(15, 23).Value = 'PANELS-(PANELS)'
$ws.Cells.Item(15, 24).Value = 'Right side panel -(Right side panel )'
$ws.Cells.Item(15, 25).Value = 'INTERIOR PANEL DIRTY & RUSTED'
$ws.Cells.Item(15, 26).Value = ""
$ws.Cells.Item(15, 26).Style = "Normal"
$ws.Cells.Item(15, 27).Value = ""
$ws.Cells.Item(15, 27).Style = "Normal"
$ws.Cells.Item(15, 28).Value = ""
$ws.Cells.Item(15, 28).Style = "Normal"

# Row 16
$ws.Cells.Item(16, 2).Value = ""
$ws.Cells.Item(16, 2).Style = "Normal"
$ws.Cells.Item(16, 3).Style = "Normal"
$ws.Cells.Item(16, 4).Style = "Normal"
$ws.Cells.Item(16, 5).Style = "Normal"
$ws.Cells.Item(16, 6).Style = "Normal"
$ws.Cells.Item(16, 7).Style = "Normal"
$ws.Cells.Item(16, 8).Style = "Normal"
$ws.Cells.Item(16, 9).Style = "Normal"
$ws.Cells.Item(16, 10).Style = "Normal"
$ws.Cells.Item(16, 11).Style = "Normal"
$ws.Cells.Item(16, 12).Style = "Normal"
$ws.Cells.Item(16, 14).Style = "Normal"
$ws.Cells.Item(16, 15).Style = "Normal"
$ws.Cells.Item(16, 16).Style = "Normal"
$ws.Cells.Item(16, 17).Style = "Normal"
$ws.Cells.Item(16, 18).Style = "Normal"
$ws.Cells.Item(16, 19).Style = "Normal"
$ws.Cells.Item(16, 20).Style = "Normal"
$ws.Cells.Item(16, 22).Style = "Normal"
$ws.Cells.Item(16, 13).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(16, 21).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(16, 23).Value = 'PANELS-(PANELS)'
$ws.Cells.Item(16, 24).Value = 'Front panel -(Front panel )'
$ws.Cells.Item(16, 25).Value = 'FRONT PANEL DENTED 36"X36"'
$ws.Cells.Item(16, 26).Value = ""
$ws.Cells.Item(16, 26).Style = "Normal"
$ws.Cells.Item(16, 27).Value = ""
$ws.Cells.Item(16, 27).Style = "Normal"
$ws.Cells.Item(16, 28).Value = ""
$ws.Cells.Item(16, 28).Style = "Normal"

# Row 22
$ws.Cells.Item(22, 2).Value = ""
$ws.Cells.Item(22, 2).Style = "Normal"
$ws.Cells.Item(22, 3).Style = "Normal"
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(22, 5).Style = "Normal"
$ws.Cells.Item(22, 6).Style = "Normal"
$ws.Cells.Item(22, 7).Style = "Normal"
$ws.Cells.Item(22, 8).Style = "Normal"
$ws.Cells.Item(22, 9).Style = "Normal"
$ws.Cells.Item(22, 10).Style = "Normal"
$ws.Cells.Item(22, 11).Style = "Normal"
$ws.Cells.Item(22, 12).Style = "Normal"
$ws.Cells.Item(22, 14).Style = "Normal"
$ws.Cells.Item(22, 15).Style = "Normal"
$ws.Cells.Item(22, 16).Style = "Normal"
$ws.Cells.Item(22, 17).Style = "Normal"
$ws.Cells.Item(22, 18).Style = "Normal"
$ws.Cells.Item(22, 19).Style = "Normal"
$ws.Cells.Item(22, 20).Style = "Normal"
$ws.Cells.Item(22, 22).Style = "Normal"
$ws.Cells.Item(22, 13).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(22, 21).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(22, 23).Value = 'PANELS-(PANELS)'
$ws.Cells.Item(22, 24).Value = 'Right side panel -(Right side panel )'
$ws.Cells.Item(22, 25).Value = 'INTERIRO PANEL DIRTY SCRATCHED'
$ws.Cells.Item(22, 26).Value = ""
$ws.Cells.Item(22, 26).Style = "Normal"
$ws.Cells.Item(22, 27).Value = ""
$ws.Cells.Item(22, 27).Style = "Normal"
$ws.Cells.Item(22, 28).Value = ""
$ws.Cells.Item(22, 28).Style = "Normal"

# Row 24
$ws.Cells.Item(24, 2).Value = ""
$ws.Cells.Item(24, 2).Style = "Normal"
$ws.Cells.Item(24, 3).Style = "Normal"
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Cells.Item(24, 5).Style = "Normal"
$ws.Cells.Item(24, 6).Style = "Normal"
$ws.Cells.Item(24, 7).Style = "Normal"
$ws.Cells.Item(24, 8).Style = "Normal"
$ws.Cells.Item(24, 9).Style = "Normal"
$ws.Cells.Item(24, 10).Style = "Normal"
$ws.Cells.Item(24, 11).Style = "Normal"
$ws.Cells.Item(24, 12).Style = "Normal"
$ws.Cells.Item(24, 14).Style = "Normal"
$ws.Cells.Item(24, 15).Style = "Normal"
$ws.Cells.Item(24, 16).Style = "Normal"
$ws.Cells.Item(24, 17).Style = "Normal"
$ws.Cells.Item(24, 18).Style = "Normal"
$ws.Cells.Item(24, 19).Style = "Normal"
$ws.Cells.Item(24, 20).Style = "Normal"
$ws.Cells.Item(24, 22).Style = "Normal"
$ws.Cells.Item(24, 13).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(24, 21).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(24, 23).Value = 'FLOORS-(F)'
$ws.Cells.Item(24, 24).Value = 'Threshold plate-(Threshold plate)'
$ws.Cells.Item(24, 25).Value = 'FLOOR BOARD WOOD FITTING 05 PCS & SCRWO MISSING 6 PCS'
$ws.Cells.Item(24, 26).Value = ""
$ws.Cells.Item(24, 26).Style = "Normal"
$ws.Cells.Item(24, 27).Value = ""
$ws.Cells.Item(24, 27).Style = "Normal"
$ws.Cells.Item(24, 28).Value = ""
$ws.Cells.Item(24, 28).Style = "Normal"

# Row 25
$ws.Cells.Item(25, 2).Value = ""
$ws.Cells.Item(25, 2).Style = "Normal"
$ws.Cells.Item(25, 3).Style = "Normal"
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(25, 5).Style = "Normal"
$ws.Cells.Item(25, 6).Style = "Normal"
$ws.Cells.Item(25, 7).Style = "Normal"
$ws.Cells.Item(25, 8).Style = "Normal"
$ws.Cells.Item(25, 9).Style = "Normal"
$ws.Cells.Item(25, 10).Style = "Normal"
$ws.Cells.Item(25, 11).Style = "Normal"
$ws.Cells.Item(25, 12).Style = "Normal"
$ws.Cells.Item(25, 14).Style = "Normal"
$ws.Cells.Item(25, 15).Style = "Normal"
$ws.Cells.Item(25, 16).Style = "Normal"
$ws.Cells.Item(25, 17).Style = "Normal"
$ws.Cells.Item(25, 18).Style = "Normal"
$ws.Cells.Item(25, 19).Style = "Normal"
$ws.Cells.Item(25, 20).Style = "Normal"
$ws.Cells.Item(25, 22).Style = "Normal"
$ws.Cells.Item(25, 13).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(25, 21).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(25, 23).Value = 'DOORS-(D)'
$ws.Cells.Item(25, 24).Value = 'Door stiffeners hinges side edge.-(Door stiffeners hinges side edge.)'
$ws.Cells.Item(25, 25).Value = 'RIGHT DOOR LOCKING BAR BENT 01 PCS'
$ws.Cells.Item(25, 26).Value = ""
$ws.Cells.Item(25, 26).Style = "Normal"
$ws.Cells.Item(25, 27).Value = ""
$ws.Cells.Item(25, 27).Style = "Normal"
$ws.Cells.Item(25, 28).Value = ""
$ws.Cells.Item(25, 28).Style = "Normal"

# Row 27
$ws.Cells.Item(27, 2).Value = ""
$ws.Cells.Item(27, 2).Style = "Normal"
$ws.Cells.Item(27, 3).Style = "Normal"
$ws.Cells.Item(27, 4).Style = "Normal"
$ws.Cells.Item(27, 5).Style = "Normal"
$ws.Cells.Item(27, 6).Style = "Normal"
$ws.Cells.Item(27, 7).Style = "Normal"
$ws.Cells.Item(27, 8).Style = "Normal"
$ws.Cells.Item(27, 9).Style = "Normal"
$ws.Cells.Item(27, 10).Style = "Normal"
$ws.Cells.Item(27, 11).Style = "Normal"
$ws.Cells.Item(27, 12).Style = "Normal"
$ws.Cells.Item(27, 14).Style = "Normal"
$ws.Cells.Item(27, 15).Style = "Normal"
$ws.Cells.Item(27, 16).Style = "Normal"
$ws.Cells.Item(27, 17).Style = "Normal"
$ws.Cells.Item(27, 18).Style = "Normal"
$ws.Cells.Item(27, 19).Style = "Normal"
$ws.Cells.Item(27, 20).Style = "Normal"
$ws.Cells.Item(27, 22).Style = "Normal"
$ws.Cells.Item(27, 13).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(27, 21).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(27, 23).Value = 'FLOORS-(F)'
$ws.Cells.Item(27, 24).Value = 'Threshold plate-(Threshold plate)'
$ws.Cells.Item(27, 25).Value = 'FLOOR BOARD WOOD FITTING 04 PEC'
$ws.Cells.Item(27, 26).Value = ""
$ws.Cells.Item(27, 26).Style = "Normal"
$ws.Cells.Item(27, 27).Value = ""
$ws.Cells.Item(27, 27).Style = "Normal"
$ws.Cells.Item(27, 28).Value = ""
$ws.Cells.Item(27, 28).Style = "Normal"

# Row 30
$ws.Cells.Item(30, 2).Value = ""
$ws.Cells.Item(30, 2).Style = "Normal"
$ws.Cells.Item(30, 3).Style = "Normal"
$ws.Cells.Item(30, 4).Style = "Normal"
$ws.Cells.Item(30, 5).Style = "Normal"
$ws.Cells.Item(30, 6).Style = "Normal"
$ws.Cells.Item(30, 7).Style = "Normal"
$ws.Cells.Item(30, 8).Style = "Normal"
$ws.Cells.Item(30, 9).Style = "Normal"
$ws.Cells.Item(30, 10).Style = "Normal"
$ws.Cells.Item(30, 11).Style = "Normal"
$ws.Cells.Item(30, 12).Style = "Normal"
$ws.Cells.Item(30, 14).Style = "Normal"
$ws.Cells.Item(30, 15).Style = "Normal"
$ws.Cells.Item(30, 16).Style = "Normal"
$ws.Cells.Item(30, 17).Style = "Normal"
$ws.Cells.Item(30, 18).Style = "Normal"
$ws.Cells.Item(30, 19).Style = "Normal"
$ws.Cells.Item(30, 20).Style = "Normal"
$ws.Cells.Item(30, 22).Style = "Normal"
$ws.Cells.Item(30, 13).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(30, 21).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(30, 23).Value = 'FLOORS-(F)'
$ws.Cells.Item(30, 24).Value = 'FLOOR BOARD-(FLOOR BOARD)'
$ws.Cells.Item(30, 25).Value = 'FLOOR BOARD BROKEN 12"X12" AT 2 PLASE,'
$ws.Cells.Item(30, 26).Value = ""
$ws.Cells.Item(30, 26).Style = "Normal"
$ws.Cells.Item(30, 27).Value = ""
$ws.Cells.Item(30, 27).Style = "Normal"
$ws.Cells.Item(30, 28).Value = ""
$ws.Cells.Item(30, 28).Style = "Normal"

# Row 31
$ws.Cells.Item(31, 2).Value = ""
$ws.Cells.Item(31, 2).Style = "Normal"
$ws.Cells.Item(31, 3).Style = "Normal"
$ws.Cells.Item(31, 4).Style = "Normal"
$ws.Cells.Item(31, 5).Style = "Normal"
$ws.Cells.Item(31, 6).Style = "Normal"
$ws.Cells.Item(31, 7).Style = "Normal"
$ws.Cells.Item(31, 8).Style = "Normal"
$ws.Cells.Item(31, 9).Style = "Normal"
$ws.Cells.Item(31, 10).Style = "Normal"
$ws.Cells.Item(31, 11).Style = "Normal"
$ws.Cells.Item(31, 12).Style = "Normal"
$ws.Cells.Item(31, 14).Style = "Normal"
$ws.Cells.Item(31, 15).Style = "Normal"
$ws.Cells.Item(31, 16).Style = "Normal"
$ws.Cells.Item(31, 17).Style = "Normal"
$ws.Cells.Item(31, 18).Style = "Normal"
$ws.Cells.Item(31, 19).Style = "Normal"
$ws.Cells.Item(31, 20).Style = "Normal"
$ws.Cells.Item(31, 22).Style = "Normal"
$ws.Cells.Item(31, 13).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(31, 21).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(31, 23).Value = 'PANELS-(PANELS)'
$ws.Cells.Item(31, 24).Value = 'Roof panel-(Roof panel)'
$ws.Cells.Item(31, 25).Value = 'ROOF PANEL CUT 8"'
$ws.Cells.Item(31, 26).Value = ""
$ws.Cells.Item(31, 26).Style = "Normal"
$ws.Cells.Item(31, 27).Value = ""
$ws.Cells.Item(31, 27).Style = "Normal"
$ws.Cells.Item(31, 28).Value = ""
$ws.Cells.Item(31, 28).Style = "Normal"

# Row 32
$ws.Cells.Item(32, 2).Value = ""
$ws.Cells.Item(32, 2).Style = "Normal"
$ws.Cells.Item(32, 3).Style = "Normal"
$ws.Cells.Item(32, 4).Style = "Normal"
$ws.Cells.Item(32, 5).Style = "Normal"
$ws.Cells.Item(32, 6).Style = "Normal"
$ws.Cells.Item(32, 7).Style = "Normal"
$ws.Cells.Item(32, 8).Style = "Normal"
$ws.Cells.Item(32, 9).Style = "Normal"
$ws.Cells.Item(32, 10).Style = "Normal"
$ws.Cells.Item(32, 11).Style = "Normal"
$ws.Cells.Item(32, 12).Style = "Normal"
$ws.Cells.Item(32, 14).Style = "Normal"
$ws.Cells.Item(32, 15).Style = "Normal"
$ws.Cells.Item(32, 16).Style = "Normal"
$ws.Cells.Item(32, 17).Style = "Normal"
$ws.Cells.Item(32, 18).Style = "Normal"
$ws.Cells.Item(32, 19).Style = "Normal"
$ws.Cells.Item(32, 20).Style = "Normal"
$ws.Cells.Item(32, 22).Style = "Normal"
$ws.Cells.Item(32, 13).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(32, 21).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(32, 23).Value = 'PANELS-(PANELS)'
$ws.Cells.Item(32, 24).Value = 'Right side panel -(Right side panel )'
$ws.Cells.Item(32, 25).Value = 'RIGHT SIDE PANEL CUT 4" AT 3 PLASE,'
$ws.Cells.Item(32, 26).Value = ""
$ws.Cells.Item(32, 26).Style = "Normal"
$ws.Cells.Item(32, 27).Value = ""
$ws.Cells.Item(32, 27).Style = "Normal"
$ws.Cells.Item(32, 28).Value = ""
$ws.Cells.Item(32, 28).Style = "Normal"

# Row 33
$ws.Cells.Item(33, 2).Value = ""
$ws.Cells.Item(33, 2).Style = "Normal"
$ws.Cells.Item(33, 3).Style = "Normal"
$ws.Cells.Item(33, 4).Style = "Normal"
$ws.Cells.Item(33, 5).Style = "Normal"
$ws.Cells.Item(33, 6).Style = "Normal"
$ws.Cells.Item(33, 7).Style = "Normal"
$ws.Cells.Item(33, 8).Style = "Normal"
$ws.Cells.Item(33, 9).Style = "Normal"
$ws.Cells.Item(33, 10).Style = "Normal"
$ws.Cells.Item(33, 11).Style = "Normal"
$ws.Cells.Item(33, 12).Style = "Normal"
$ws.Cells.Item(33, 14).Style = "Normal"
$ws.Cells.Item(33, 15).Style = "Normal"
$ws.Cells.Item(33, 16).Style = "Normal"
$ws.Cells.Item(33, 17).Style = "Normal"
$ws.Cells.Item(33, 18).Style = "Normal"
$ws.Cells.Item(33, 19).Style = "Normal"
$ws.Cells.Item(33, 20).Style = "Normal"
$ws.Cells.Item(33, 22).Style = "Normal"
$ws.Cells.Item(33, 13).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(33, 21).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(33, 23).Value = 'DOORS-(D)'
$ws.Cells.Item(33, 24).Value = 'Door stiffeners hinges side edge.-(Door stiffeners hinges side edge.)'
$ws.Cells.Item(33, 25).Value = 'BOTH DOOR HINGS BENT 4 PCS(DOOR UN LOCK)'
$ws.Cells.Item(33, 26).Value = ""
$ws.Cells.Item(33, 26).Style = "Normal"
$ws.Cells.Item(33, 27).Value = ""
$ws.Cells.Item(33, 27).Style = "Normal"
$ws.Cells.Item(33, 28).Value = ""
$ws.Cells.Item(33, 28).Style = "Normal"

# Row 35
$ws.Cells.Item(35, 2).Value = ""
$ws.Cells.Item(35, 2).Style = "Normal"
$ws.Cells.Item(35, 3).Style = "Normal"
$ws.Cells.Item(35, 4).Style = "Normal"
$ws.Cells.Item(35, 5).Style = "Normal"
$ws.Cells.Item(35, 6).Style = "Normal"
$ws.Cells.Item(35, 7).Style = "Normal"
$ws.Cells.Item(35, 8).Style = "Normal"
$ws.Cells.Item(35, 9).Style = "Normal"
$ws.Cells.Item(35, 10).Style = "Normal"
$ws.Cells.Item(35, 11).Style = "Normal"
$ws.Cells.Item(35, 12).Style = "Normal"
$ws.Cells.Item(35, 14).Style = "Normal"
$ws.Cells.Item(35, 15).Style = "Normal"
$ws.Cells.Item(35, 16).Style = "Normal"
$ws.Cells.Item(35, 17).Style = "Normal"
$ws.Cells.Item(35, 18).Style = "Normal"
$ws.Cells.Item(35, 19).Style = "Normal"
$ws.Cells.Item(35, 20).Style = "Normal"
$ws.Cells.Item(35, 22).Style = "Normal"
$ws.Cells.Item(35, 13).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(35, 21).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(35, 23).Value = 'FLOORS-(F)'
$ws.Cells.Item(35, 24).Value = 'FLOOR BOARD-(FLOOR BOARD)'
$ws.Cells.Item(35, 25).Value = 'FLOOR BOARD PICKUP 24"X24"'
$ws.Cells.Item(35, 26).Value = ""
$ws.Cells.Item(35, 26).Style = "Normal"
$ws.Cells.Item(35, 27).Value = ""
$ws.Cells.Item(35, 27).Style = "Normal"
$ws.Cells.Item(35, 28).Value = ""
$ws.Cells.Item(35, 28).Style = "Normal"

# Row 36
$ws.Cells.Item(36, 2).Value = ""
$ws.Cells.Item(36, 2).Style = "Normal"
$ws.Cells.Item(36, 3).Style = "Normal"
$ws.Cells.Item(36, 4).Style = "Normal"
$ws.Cells.Item(36, 5).Style = "Normal"
$ws.Cells.Item(36, 6).Style = "Normal"
$ws.Cells.Item(36, 7).Style = "Normal"
$ws.Cells.Item(36, 8).Style = "Normal"
$ws.Cells.Item(36, 9).Style = "Normal"
$ws.Cells.Item(36, 10).Style = "Normal"
$ws.Cells.Item(36, 11).Style = "Normal"
$ws.Cells.Item(36, 12).Style = "Normal"
$ws.Cells.Item(36, 14).Style = "Normal"
$ws.Cells.Item(36, 15).Style = "Normal"
$ws.Cells.Item(36, 16).Style = "Normal"
$ws.Cells.Item(36, 17).Style = "Normal"
$ws.Cells.Item(36, 18).Style = "Normal"
$ws.Cells.Item(36, 19).Style = "Normal"
$ws.Cells.Item(36, 20).Style = "Normal"
$ws.Cells.Item(36, 22).Style = "Normal"
$ws.Cells.Item(36, 13).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(36, 21).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(36, 23).Value = 'PANELS-(PANELS)'
$ws.Cells.Item(36, 24).Value = 'Front panel -(Front panel )'
$ws.Cells.Item(36, 25).Value = 'FRONT PANEL CUT 2"'
$ws.Cells.Item(36, 26).Value = ""
$ws.Cells.Item(36, 26).Style = "Normal"
$ws.Cells.Item(36, 27).Value = ""
$ws.Cells.Item(36, 27).Style = "Normal"
$ws.Cells.Item(36, 28).Value = ""
$ws.Cells.Item(36, 28).Style = "Normal"

# Row 37
$ws.Cells.Item(37, 2).Value = ""
$ws.Cells.Item(37, 2).Style = "Normal"
$ws.Cells.Item(37, 3).Style = "Normal"
$ws.Cells.Item(37, 4).Style = "Normal"
$ws.Cells.Item(37, 5).Style = "Normal"
$ws.Cells.Item(37, 6).Style = "Normal"
$ws.Cells.Item(37, 7).Style = "Normal"
$ws.Cells.Item(37, 8).Style = "Normal"
$ws.Cells.Item(37, 9).Style = "Normal"
$ws.Cells.Item(37, 10).Style = "Normal"
$ws.Cells.Item(37, 11).Style = "Normal"
$ws.Cells.Item(37, 12).Style = "Normal"
$ws.Cells.Item(37, 14).Style = "Normal"
$ws.Cells.Item(37, 15).Style = "Normal"
$ws.Cells.Item(37, 16).Style = "Normal"
$ws.Cells.Item(37, 17).Style = "Normal"
$ws.Cells.Item(37, 18).Style = "Normal"
$ws.Cells.Item(37, 19).Style = "Normal"
$ws.Cells.Item(37, 20).Style = "Normal"
$ws.Cells.Item(37, 22).Style = "Normal"
$ws.Cells.Item(37, 13).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(37, 21).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(37, 23).Value = 'PANELS-(PANELS)'
$ws.Cells.Item(37, 24).Value = 'Right side panel -(Right side panel )'
$ws.Cells.Item(37, 25).Value = 'RIGHT SIDE PANEL CUT 2"X2" AT 4 PLASE,'
$ws.Cells.Item(37, 26).Value = ""
$ws.Cells.Item(37, 26).Style = "Normal"
$ws.Cells.Item(37, 27).Value = ""
$ws.Cells.Item(37, 27).Style = "Normal"
$ws.Cells.Item(37, 28).Value = ""
$ws.Cells.Item(37, 28).Style = "Normal"

# Row 38
$ws.Cells.Item(38, 2).Value = ""
$ws.Cells.Item(38, 2).Style = "Normal"
$ws.Cells.Item(38, 3).Style = "Normal"
$ws.Cells.Item(38, 4).Style = "Normal"
$ws.Cells.Item(38, 5).Style = "Normal"
$ws.Cells.Item(38, 6).Style = "Normal"
$ws.Cells.Item(38, 7).Style = "Normal"
$ws.Cells.Item(38, 8).Style = "Normal"
$ws.Cells.Item(38, 9).Style = "Normal"
$ws.Cells.Item(38, 10).Style = "Normal"
$ws.Cells.Item(38, 11).Style = "Normal"
$ws.Cells.Item(38, 12).Style = "Normal"
$ws.Cells.Item(38, 14).Style = "Normal"
$ws.Cells.Item(38, 15).Style = "Normal"
$ws.Cells.Item(38, 16).Style = "Normal"
$ws.Cells.Item(38, 17).Style = "Normal"
$ws.Cells.Item(38, 18).Style = "Normal"
$ws.Cells.Item(38, 19).Style = "Normal"
$ws.Cells.Item(38, 20).Style = "Normal"
$ws.Cells.Item(38, 22).Style = "Normal"
$ws.Cells.Item(38, 13).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(38, 21).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(38, 23).Value = 'PANELS-(PANELS)'
$ws.Cells.Item(38, 24).Value = 'Left side panel-(Left side panel)'
$ws.Cells.Item(38, 25).Value = 'LEFT SIDE PANEL CUT 2"X2"'
$ws.Cells.Item(38, 26).Value = ""
$ws.Cells.Item(38, 26).Style = "Normal"
$ws.Cells.Item(38, 27).Value = ""
$ws.Cells.Item(38, 27).Style = "Normal"
$ws.Cells.Item(38, 28).Value = ""
$ws.Cells.Item(38, 28).Style = "Normal"

# Row 39
$ws.Cells.Item(39, 2).Value = ""
$ws.Cells.Item(39, 2).Style = "Normal"
$ws.Cells.Item(39, 3).Style = "Normal"
$ws.Cells.Item(39, 4).Style = "Normal"
$ws.Cells.Item(39, 5).Style = "Normal"
$ws.Cells.Item(39, 6).Style = "Normal"
$ws.Cells.Item(39, 7).Style = "Normal"
$ws.Cells.Item(39, 8).Style = "Normal"
$ws.Cells.Item(39, 9).Style = "Normal"
$ws.Cells.Item(39, 10).Style = "Normal"
$ws.Cells.Item(39, 11).Style = "Normal"
$ws.Cells.Item(39, 12).Style = "Normal"
$ws.Cells.Item(39, 14).Style = "Normal"
$ws.Cells.Item(39, 15).Style = "Normal"
$ws.Cells.Item(39, 16).Style = "Normal"
$ws.Cells.Item(39, 17).Style = "Normal"
$ws.Cells.Item(39, 18).Style = "Normal"
$ws.Cells.Item(39, 19).Style = "Normal"
$ws.Cells.Item(39, 20).Style = "Normal"
$ws.Cells.Item(39, 22).Style = "Normal"
$ws.Cells.Item(39, 13).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(39, 21).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(39, 23).Value = 'PANELS-(PANELS)'
$ws.Cells.Item(39, 24).Value = 'Right side panel -(Right side panel )'
$ws.Cells.Item(39, 25).Value = 'BOTH SIDE PANEL PUSH OUT 40"X40"'
$ws.Cells.Item(39, 26).Value = ""
$ws.Cells.Item(39, 26).Style = "Normal"
$ws.Cells.Item(39, 27).Value = ""
$ws.Cells.Item(39, 27).Style = "Normal"
$ws.Cells.Item(39, 28).Value = ""
$ws.Cells.Item(39, 28).Style = "Normal"

# Row 43
$ws.Cells.Item(43, 2).Value = ""
$ws.Cells.Item(43, 2).Style = "Normal"
$ws.Cells.Item(43, 3).Style = "Normal"
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Cells.Item(43, 5).Style = "Normal"
$ws.Cells.Item(43, 6).Style = "Normal"
$ws.Cells.Item(43, 7).Style = "Normal"
$ws.Cells.Item(43, 8).Style = "Normal"
$ws.Cells.Item(43, 9).Style = "Normal"
$ws.Cells.Item(43, 10).Style = "Normal"
$ws.Cells.Item(43, 11).Style = "Normal"
$ws.Cells.Item(43, 12).Style = "Normal"
$ws.Cells.Item(43, 14).Style = "Normal"
$ws.Cells.Item(43, 15).Style = "Normal"
$ws.Cells.Item(43, 16).Style = "Normal"
$ws.Cells.Item(43, 17).Style = "Normal"
$ws.Cells.Item(43, 18).Style = "Normal"
$ws.Cells.Item(43, 19).Style = "Normal"
$ws.Cells.Item(43, 20).Style = "Normal"
$ws.Cells.Item(43, 22).Style = "Normal"
$ws.Cells.Item(43, 13).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(43, 21).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(43, 23).Value = 'PANELS-(PANELS)'
$ws.Cells.Item(43, 24).Value = 'Right side panel -(Right side panel )'
$ws.Cells.Item(43, 25).Value = 'RIGHT SIDE PANEL CUT 2" & LEFT SIDE PANEL CUT 2" AT  2 PLC '
$ws.Cells.Item(43, 26).Value = ""
$ws.Cells.Item(43, 26).Style = "Normal"
$ws.Cells.Item(43, 27).Value = ""
$ws.Cells.Item(43, 27).Style = "Normal"
$ws.Cells.Item(43, 28).Value = ""
$ws.Cells.Item(43, 28).Style = "Normal"

# Row 44
$ws.Cells.Item(44, 2).Value = ""
$ws.Cells.Item(44, 2).Style = "Normal"
$ws.Cells.Item(44, 3).Style = "Normal"
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Cells.Item(44, 5).Style = "Normal"
$ws.Cells.Item(44, 6).Style = "Normal"
$ws.Cells.Item(44, 7).Style = "Normal"
$ws.Cells.Item(44, 8).Style = "Normal"
$ws.Cells.Item(44, 9).Style = "Normal"
$ws.Cells.Item(44, 10).Style = "Normal"
$ws.Cells.Item(44, 11).Style = "Normal"
$ws.Cells.Item(44, 12).Style = "Normal"
$ws.Cells.Item(44, 14).Style = "Normal"
$ws.Cells.Item(44, 15).Style = "Normal"
$ws.Cells.Item(44, 16).Style = "Normal"
$ws.Cells.Item(44, 17).Style = "Normal"
$ws.Cells.Item(44, 18).Style = "Normal"
$ws.Cells.Item(44, 19).Style = "Normal"
$ws.Cells.Item(44, 20).Style = "Normal"
$ws.Cells.Item(44, 22).Style = "Normal"
$ws.Cells.Item(44, 13).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(44, 21).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(44, 23).Value = 'PANELS-(PANELS)'
$ws.Cells.Item(44, 24).Value = 'Right side panel -(Right side panel )'
$ws.Cells.Item(44, 25).Value = 'BOTH OUT SIDE PANEL TWISTED '
$ws.Cells.Item(44, 26).Value = ""
$ws.Cells.Item(44, 26).Style = "Normal"
$ws.Cells.Item(44, 27).Value = ""
$ws.Cells.Item(44, 27).Style = "Normal"
$ws.Cells.Item(44, 28).Value = ""
$ws.Cells.Item(44, 28).Style = "Normal"

# Row 46
$ws.Cells.Item(46, 2).Value = ""
$ws.Cells.Item(46, 2).Style = "Normal"
$ws.Cells.Item(46, 3).Style = "Normal"
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Cells.Item(46, 5).Style = "Normal"
$ws.Cells.Item(46, 6).Style = "Normal"
$ws.Cells.Item(46, 7).Style = "Normal"
$ws.Cells.Item(46, 8).Style = "Normal"
$ws.Cells.Item(46, 9).Style = "Normal"
$ws.Cells.Item(46, 10).Style = "Normal"
$ws.Cells.Item(46, 11).Style = "Normal"
$ws.Cells.Item(46, 12).Style = "Normal"
$ws.Cells.Item(46, 14).Style = "Normal"
$ws.Cells.Item(46, 15).Style = "Normal"
$ws.Cells.Item(46, 16).Style = "Normal"
$ws.Cells.Item(46, 17).Style = "Normal"
$ws.Cells.Item(46, 18).Style = "Normal"
$ws.Cells.Item(46, 19).Style = "Normal"
$ws.Cells.Item(46, 20).Style = "Normal"
$ws.Cells.Item(46, 22).Style = "Normal"
$ws.Cells.Item(46, 13).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(46, 21).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(46, 23).Value = 'DOORS-(D)'
$ws.Cells.Item(46, 24).Value = 'Door stiffeners hinges side edge.-(Door stiffeners hinges side edge.)'
$ws.Cells.Item(46, 25).Value = 'RIGHT DOOR BOTTOM GASKET CUT 8"'
$ws.Cells.Item(46, 26).Value = ""
$ws.Cells.Item(46, 26).Style = "Normal"
$ws.Cells.Item(46, 27).Value = ""
$ws.Cells.Item(46, 27).Style = "Normal"
$ws.Cells.Item(46, 28).Value = ""
$ws.Cells.Item(46, 28).Style = "Normal"

